$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry over the header row's formatting (bold, centered/top, thin border)
# to the new row's first cell, exactly like Excel's "copy formatting" does.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Populate the new row of data (row 2).
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "sss4"
$ws.Range("C2").Value = "sss4"
$ws.Range("D2").Value = "sss5"
$ws.Range("E2").Value = "sss6"
$ws.Range("F2").Value = "sss8"
$ws.Range("G2").Value = "sss10"
